$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.455.74'
$ws.Range("E2").Value = '  +5.82%  '
$ws.Range("D3").Value = '2.514.98'
$ws.Range("E3").Value = '  +3.85%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''324.61'
$ws.Range("E5").Value = '  +2.29%  '
$ws.Range("D6").Value = '''105.66'
$ws.Range("E6").Value = '  +2.45%  '
$ws.Range("E7").Value = '  +1.96%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +3.25%  '
$ws.Range("D10").Value = '''37.23'
$ws.Range("E10").Value = '  +4.91%  '
$ws.Range("D11").Value = '''0.0820'
$ws.Range("E11").Value = '  +2.35%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '''18.56'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = '''7.22'
$ws.Range("E14").Value = '  +4.53%  '
$ws.Range("D15").Value = '2.911.24'
$ws.Range("E15").Value = '  +4.05%  '
$ws.Range("D16").Value = '2.525.86'
$ws.Range("E16").Value = '  +4.65%  '
$ws.Range("E17").Value = '  +2.84%  '
$ws.Range("D18").Value = '47.410.67'
$ws.Range("E18").Value = '  +6.08%  '
$ws.Range("D19").Value = '''12.83'
$ws.Range("E19").Value = '  +4.43%  '
$ws.Range("D20").Value = '''6.62'
$ws.Range("E20").Value = '  +4.33%  '
$ws.Range("E21").Value = '  +2.93%  '
$ws.Range("D22").Value = '''71.18'
$ws.Range("D23").Value = '''252.73'
$ws.Range("E23").Value = '  +3.87%  '
$ws.Range("E24").Value = '  +5.07%  '
$ws.Range("E25").Value = '  +3.25%  '
$ws.Range("D26").Value = '''26.50'
$ws.Range("E26").Value = '  +5.05%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = '''10.07'
$ws.Range("E28").Value = '  +5.53%  '
$ws.Range("D29").Value = '''2.20'
$ws.Range("E29").Value = '  -3.60%  '
$ws.Range("D30").Value = '''35.37'
$ws.Range("E30").Value = '  +5.59%  '
$ws.Range("D31").Value = '''0.136'
$ws.Range("E31").Value = '  +7.13%  '
$ws.Range("D32").Value = '''49.66'
$ws.Range("E32").Value = '  +2.31%  '
$ws.Range("D33").Value = '''19.88'
$ws.Range("E33").Value = '  +1.80%  '
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("E35").Value = '  +3.16%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = '''1.95'
$ws.Range("E37").Value = '  +3.70%  '
$ws.Range("E38").Value = '  +4.19%  '
$ws.Range("E39").Value = '  +4.80%  '
$ws.Range("D40").Value = '''123.59'
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("E41").Value = '  +2.38%  '
$ws.Range("D42").Value = '''2.25'
$ws.Range("E42").Value = '  +3.30%  '
$ws.Range("D43").Value = '''21.72'
$ws.Range("E43").Value = '  +3.87%  '
$ws.Range("D44").Value = '''0.0299'
$ws.Range("E44").Value = '  +3.72%  '
$ws.Range("D45").Value = '1.984.82'
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("D46").Value = '''3.03'
$ws.Range("E46").Value = '  +3.45%  '
$ws.Range("D47").Value = '''2.13'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").Value = '''1.80'
$ws.Range("E48").Value = '  +2.31%  '
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").Value = '''5.46'
$ws.Range("E50").Value = '  +18.09%  '
$ws.Range("D51").Value = '''79.78'
$ws.Range("E51").Value = '  +4.94%  '
